$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.784.67"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "1.633.77"
$ws.Range("E3").Value = "  +0.23%  "
$__style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = $__style
$ws.Range("E4").Value = "  +0.21%  "
$__style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.95"
$ws.Range("D5").Style = $__style
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("E6").Value = "  -0.74%  "
$__style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("D7").Style = $__style
$ws.Range("E7").Value = "  +0.19%  "
$__style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2576"
$ws.Range("D8").Style = $__style
$ws.Range("E8").Value = "  +0.68%  "
$__style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06419"
$ws.Range("D9").Style = $__style
$ws.Range("E9").Value = "  +1.40%  "
$__style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.21"
$ws.Range("D10").Style = $__style
$ws.Range("E10").Value = "  +3.95%  "
$__style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07787"
$ws.Range("D11").Style = $__style
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.654.78"
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$__style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.248"
$ws.Range("D13").Style = $__style
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "1.861.77"
$ws.Range("E14").Value = "  +0.35%  "
$__style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5576"
$ws.Range("D15").Style = $__style
$ws.Range("E15").Value = "  +1.06%  "
$ws.Range("D16").Value = "0.0₅7648"
$ws.Range("E16").Value = "  +0.62%  "
$__style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.18"
$ws.Range("D17").Style = $__style
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").Value = "25.807.61"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("E19").Value = "  +0.19%  "
$__style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.363"
$ws.Range("D20").Style = $__style
$__style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "191.90"
$ws.Range("D21").Style = $__style
$ws.Range("E21").Value = "  -1.37%  "
$__style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.899"
$ws.Range("D22").Style = $__style
$ws.Range("E22").Value = "  +0.46%  "
$__style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.134"
$ws.Range("D23").Style = $__style
$ws.Range("E23").Value = "  +2.00%  "
$ws.Range("E24").Value = "  +0.12%  "
$__style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.772"
$ws.Range("D25").Style = $__style
$ws.Range("E25").Value = "  -6.29%  "
$__style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "138.89"
$ws.Range("D26").Style = $__style
$ws.Range("E26").Value = "  -2.15%  "
$__style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1227"
$ws.Range("D27").Style = $__style
$__style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.809"
$ws.Range("D28").Style = $__style
$ws.Range("E28").Value = "  +0.82%  "
$ws.Range("E29").Value = "  -0.06%  "
$__style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.241"
$ws.Range("D30").Style = $__style
$ws.Range("E30").Value = "  +0.14%  "
$__style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04943"
$ws.Range("D31").Style = $__style
$ws.Range("E31").Value = "  +1.03%  "
$ws.Range("E32").Value = "  +1.87%  "
$__style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.250"
$ws.Range("D33").Style = $__style
$ws.Range("E33").Value = "  +2.31%  "
$__style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.569"
$ws.Range("D34").Style = $__style
$ws.Range("E34").Value = "  +1.53%  "
$__style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.385"
$ws.Range("D35").Style = $__style
$ws.Range("E35").Value = "  +0.51%  "
$__style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9015"
$ws.Range("D36").Style = $__style
$ws.Range("E36").Value = "  +0.83%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$__style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5574"
$ws.Range("D37").Style = $__style
$ws.Range("E37").Value = "  +1.16%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$__style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.567"
$ws.Range("D38").Style = $__style
$ws.Range("E38").Value = "  +1.19%  "
$ws.Range("D39").Value = "1.129.79"
$ws.Range("E39").Value = "  +1.39%  "
$__style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01567"
$ws.Range("D40").Style = $__style
$ws.Range("E40").Value = "  +0.95%  "
$__style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9955"
$ws.Range("D41").Style = $__style
$ws.Range("E41").Value = "  -0.47%  "
$__style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.450"
$ws.Range("D42").Style = $__style
$ws.Range("E42").Value = "  -2.19%  "
$__style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "98.86"
$ws.Range("D43").Style = $__style
$ws.Range("E43").Value = "  +1.15%  "
$__style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7983"
$ws.Range("D44").Style = $__style
$ws.Range("E44").Value = "  +0.59%  "
$ws.Range("D45").Value = "0.0₈114"
$ws.Range("E45").Value = "  +0.14%  "
$__style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.60"
$ws.Range("D46").Style = $__style
$ws.Range("E46").Value = "  +1.66%  "
$__style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4259"
$ws.Range("D47").Style = $__style
$ws.Range("E47").Value = "  -3.90%  "
$__style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.757"
$ws.Range("D48").Style = $__style
$ws.Range("E48").Value = "  +3.12%  "
$__style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05026"
$ws.Range("D49").Style = $__style
$ws.Range("E49").Value = "  -2.01%  "
$__style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.9968"
$ws.Range("D50").Style = $__style
$ws.Range("E50").Value = "  -0.49%  "
$ws.Range("E51").Value = "  +0.24%  "
